$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.792.39'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.548.74'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.56'
$ws.Range("E5").Value = '  -2.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.69'
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.95'
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '2.941.02'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.00'
$ws.Range("E15").Value = '  +5.17%  '
$ws.Range("D16").Value = '2.546.93'
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("D18").Value = '42.823.16'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.76'
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.45'
$ws.Range("E20").Value = '  -3.29%  '
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.48'
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.41'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.71'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -2.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.17'
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.18'
$ws.Range("E30").Value = '  -2.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.77'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0803'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("E37").Value = '  +4.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.38'
$ws.Range("E38").Value = '  -3.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.50'
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.12'
$ws.Range("E42").Value = '  +6.64%  '
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.22'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '1.996.40'
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.10'
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").Value = '2.786.62'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.35'
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.193'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.63'
$ws.Range("E51").Value = '  -3.62%  '
